# Auto-generated edit script to apply symbol-list refresh update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue "D2" '331.34'
Set-TextValue "E2" '0.65%'
Set-TextValue "D3" '44.72'
Set-TextValue "E3" '1.41%'
Set-TextValue "D4" '5.542'
Set-TextValue "E4" '-0.74%'
Set-TextValue "D5" '0.08157'
Set-TextValue "E5" '1.04%'
Set-TextValue "D6" '2.057'
Set-TextValue "E6" '3.99%'
Set-TextValue "B7" 'GateToken'
Set-TextValue "C7" 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue "D7" '4.431'
Set-TextValue "E7" '2.39%'
Set-TextValue "B8" 'MXToken'
Set-TextValue "C8" 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue "D8" '0.9785'
Set-TextValue "E8" '2.69%'
Set-TextValue "B9" 'LiechtensteinCryptoassetsExchange'
Set-TextValue "C9" 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue "D9" '0.1117'
Set-TextValue "E9" '-3.65%'
Set-TextValue "B10" 'WazirX'
Set-TextValue "C10" 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue "D10" '0.1898'
Set-TextValue "E10" '2.21%'
Set-TextValue "B11" 'MCDex'
Set-TextValue "C11" 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue "D11" '10.25'
Set-TextValue "E11" '-13.66%'
Set-TextValue "B12" 'MandalaExchangeToken'
Set-TextValue "C12" 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue "D12" '0.1004'
Set-TextValue "E12" '2.88%'
Set-TextValue "B13" 'BitrueCoin'
Set-TextValue "C13" 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue "D13" '0.04714'
Set-TextValue "E13" '0.68%'
Set-TextValue "B14" 'BitMartToken'
Set-TextValue "C14" 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue "D14" '0.1058'
Set-TextValue "E14" '-1.02%'
Set-TextValue "B15" 'BitForexToken'
Set-TextValue "C15" 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue "D15" '0.001286'
Set-TextValue "E15" '-0.14%'
Set-TextValue "B16" 'CoinExToken'
Set-TextValue "C16" 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue "D16" '0.04092'
Set-TextValue "E16" '-3.38%'
Set-TextValue "B17" 'TigerCash'
Set-TextValue "C17" 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue "D17" '0.005902'
Set-TextValue "E17" '-0.81%'
Set-TextValue "B18" 'HotbitToken'
Set-TextValue "C18" 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
Set-TextValue "D18" '0.004417'
Set-TextValue "E18" '1.84%'
Set-TextValue "B19" 'LEO'
Set-TextValue "C19" 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue "D19" '3.349'
Set-TextValue "E19" '-0.72%'
Set-TextValue "B20" 'BTSEToken'
Set-TextValue "C20" 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue "D20" '2.645'
Set-TextValue "E20" '2.92%'
Set-TextValue "B21" 'BitpandaEcosystemToken'
Set-TextValue "C21" 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
Set-TextValue "D21" '0.3351'
Set-TextValue "E21" '-3.55%'
Set-TextValue "B22" 'ProBitToken'
Set-TextValue "C22" 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
Set-TextValue "D22" '0.1390'
Set-TextValue "E22" '-1.48%'
Set-TextValue "B23" 'ZBToken'
Set-TextValue "C23" 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
Set-TextValue "D23" '0.2566'
Set-TextValue "E23" '2.36%'
Set-TextValue "B24" 'BitKan'
Set-TextValue "C24" 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
Set-TextValue "D24" '0.001299'
Set-TextValue "E24" '3.52%'
Set-TextValue "D25" '0.0001278'
Set-TextValue "E25" '7.30%'
Set-TextValue "D26" '0.0003733'
Set-TextValue "E26" '-6.18%'
Set-TextValue "D38" '0.02734'
Set-TextValue "E38" '3.69%'
Set-TextValue "D39" '0.05724'
Set-TextValue "E39" '3.28%'
Set-TextValue "D40" '0.007589'
Set-TextValue "E40" '0.07%'
Set-TextValue "D41" '0.1426'
Set-TextValue "E41" '1.26%'
Set-TextValue "D42" '0.007544'
Set-TextValue "E42" '-6.73%'
Set-TextValue "D43" '0.001955'
Set-TextValue "E43" '-3.08%'
Set-TextValue "D44" '0.008314'
Set-TextValue "E44" '-6.64%'
Set-TextValue "D45" '0.00007048'
Set-TextValue "E45" '-2.37%'
Set-TextValue "D46" '0.00000000749'
Set-TextValue "E46" '-0.34%'
Set-TextValue "D47" '0.0005789'
Set-TextValue "E47" '-0.38%'
Set-TextValue "D48" '0.002515'
Set-TextValue "E48" '10.64%'
Set-TextValue "D49" '0.003535'
Set-TextValue "E49" '-27.01%'
Set-TextValue "D50" '0.00002096'
Set-TextValue "E50" '-0.34%'
Set-TextValue "D51" '0.0001996'
Set-TextValue "E51" '-0.34%'

Write-Output "Applied crypto symbol list update"
